# Update the "想去人数" (number of interested attendees) counts in the
# 展览 (sheet1) and 全部类型 (sheet4) worksheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 93
    $ws.Range("F4").Value = 59
    $ws.Range("F5").Value = 2449
}
